$d = $word.ActiveDocument

# --- Title ---
$d.Content.Find.Execute("Resonance and Harmony: A Fusion of Art and Science", $true, $false, $false, $false, $false, $true, 1, $false, "Unveiling the Enigmatic Tapestry of Our Universe: A Journey Through Science and Art", 2) | Out-Null

# --- Author ---
$d.Content.Find.Execute("David Armstrong", $true, $false, $false, $false, $false, $true, 1, $false, "Evelyn Walker", 2) | Out-Null

# --- Email (two runs split by a literal "." run) ---
$d.Content.Find.Execute("musicandphysics@imaginativefusion", $true, $false, $false, $false, $false, $true, 1, $false, "ewalker2@institution", 2) | Out-Null
$d.Content.Find.Execute("net", $true, $false, $false, $false, $false, $true, 1, $false, "edu", 2) | Out-Null

# --- Body paragraph (paragraph 5): replace wholesale with the new narrative, preserving per-sentence run/rPr structure ---
$body = $d.Paragraphs(5).Range
$bodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>In the vast expanse of our cosmos, where celestial bodies dance in harmonious chaos and nature's mysteries unfold with each passing moment, we find ourselves captivated by the intricate interplay of science and art</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> From the delicate petals of a flower to the mind-boggling vastness of the universe, these seemingly disparate realms intertwine to form a unified tapestry of knowledge and beauty</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>As we delve into the depths of scientific inquiry, we uncover the fundamental principles that govern the intricate machinations of our world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Mathematics, the universal language of science, provides a framework for understanding the patterns and relationships that underpin the natural world</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Through the meticulous observation and experimentation characteristic of chemistry and biology, we unravel the secrets of life and the intricate processes that sustain it</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>Simultaneously, art serves as a conduit for expressing the ineffable, capturing the essence of human experience in myriad forms</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Whether it be through the strokes of a paintbrush, the melody of a symphony, or the written word, art enables us to transcend the limitations of language and communicate emotions, ideas, and narratives that resonate deeply within our souls</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>Introduction Continued:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>From the breathtaking grandeur of the cosmos to the microscopic marvels of cellular life, science and art unveil the wonders of our universe in complementary ways</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> The artistic spirit finds inspiration in the complexities of nature, translating scientific discoveries into captivating expressions of awe and wonder</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> Conversely, science draws upon the insights of art to illuminate the hidden patterns and symmetries that underlie the fabric of reality</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
        <w:t>Introduction Continued:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:br/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:br/>
        <w:t>As we embark on this intellectual odyssey, we gain a profound appreciation for the profound interconnectedness of science and art</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> They are not merely separate disciplines but rather two sides of the same coin, each enriching and informing the other</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve"> By embracing both rational inquiry and creative expression, we unlock a deeper understanding of the universe and our place within it</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$body.InsertXML($bodyXml)

# --- Summary paragraph (last paragraph of the "Summary" section) ---
$summary = $d.Paragraphs(7).Range
$summaryXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>In this essay, we have explored the fascinating interplay between science and art, revealing their profound interconnectedness</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> From the mathematical underpinnings of the cosmos to the artistry of nature's designs, we have seen how these fields converge to create a comprehensive understanding of our universe</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve"> This convergence inspires us to embrace both rational inquiry and creative expression as we continue to unravel the mysteries that surround us</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
'@
$summary.InsertXML($summaryXml)

# --- Trailing empty paragraph at the end of the document ---
$endRange = $d.Paragraphs(7).Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
